# Update column D (algorithm result) values for the RandomForest imputation
# result sheet. Commit message: "Update Name of Algo"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column D
$updates = @{
    11 = -8.051999999999998
    12 = -6.834499999999997
    15 = -8.594000000000001
    27 = -8.732100000000004
    28 = -8.259999999999998
    31 = -8.537799999999997
    32 = -8.9186
    36 = -8.092399999999998
    38 = -8.691899999999999
    46 = -8.357699999999999
    54 = -8.256600000000006
    55 = -8.332399999999996
    56 = -7.949299999999996
    67 = -6.406899999999999
    69 = -7.253199999999996
    72 = -7.365899999999999
    73 = -7.751299999999999
    83 = -8.605500000000001
    86 = -7.605599999999994
    91 = -6.641499999999999
    93 = -6.373999999999999
    99 = -7.699100000000004
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
